$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36/37 swap: name/link columns (B, C) ---
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"

# --- Plain text-safe updates (percentages + multi-dot / unicode prices) ---
$ws.Range("D2").Value = "62.957.00"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.547.68"
$ws.Range("E3").Value = "  +3.28%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("D9").Value = "2.546.69"
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("D15").Value = "2.999.85"
$ws.Range("E15").Value = "  +3.09%  "
$ws.Range("D16").Value = "62.931.03"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "2.541.94"
$ws.Range("E18").Value = "  +2.90%  "
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  -1.22%  "
$ws.Range("E25").Value = "  -3.53%  "
$ws.Range("E26").Value = "  +5.75%  "
$ws.Range("E27").Value = "  +11.64%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("E29").Value = "  +2.52%  "
$ws.Range("E30").Value = "  +6.72%  "
$ws.Range("D31").Value = "0.0₃0808"
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  +4.62%  "
$ws.Range("E35").Value = "  +9.52%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("E43").Value = "  +0.94%  "
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  +3.95%  "
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("E51").Value = "  -0.61%  "

# --- Numeric-looking price updates: force text via NumberFormat, then restore default style ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("D13").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "334.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "406.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.399"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "152.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.76"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.606"
$ws.Range("D46").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.26"
$ws.Range("D50").Style = "Normal"
